# Generate Report for Handoff
# Update the localization-status workbook so every reference to the old
# handoff GUID / content-hash is replaced with the new ones, and the
# recorded handoff timestamps are updated.

$wb = $excel.ActiveWorkbook

$oldGuid = "82646639-7a00-4e04-8738-375c48e838c4"
$newGuid = "583782f7-bd86-418c-a1bf-d974d024d95d"
$oldHash = "fb4f517054801433fa0dfbaab87a9b6d86ca05e4"
$newHash = "52b53c996e8dd9fae1e53aa853537a05f2c44543"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value2 = "$newGuid.md"
$ws1.Range("D2").Value2 = "2016-49-20 04:49:20"

$addrMd1 = "https://github.com/OpenLocalizationTest/oltest/blob/589670b0cadae2406c8f632e7512a04ff377a70b/e2e/$oldGuid.md"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $addrMd1, "", "", "$newGuid.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value2 = "$newGuid.md"
$ws2.Range("D2").Value2 = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("E2").Value2 = "2016-03-20 04:49:17"

$addrMd2 = "https://github.com/OpenLocalizationTest/oltest/blob/589670b0cadae2406c8f632e7512a04ff377a70b/e2e/$oldGuid.md"
$addrXlf2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0683f82c83a328f971f682780993bfeba9676c37/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $addrMd2, "", "", "$newGuid.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), $addrMd2, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), $addrXlf2, "", "", "$newGuid.$newHash.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value2 = "$newGuid.md"
$ws3.Range("D2").Value2 = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("E2").Value2 = "2016-03-20 04:49:20"

$addrMd3 = "https://github.com/OpenLocalizationTest/oltest/blob/589670b0cadae2406c8f632e7512a04ff377a70b/e2e/$oldGuid.md"
$addrXlf3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8851b3581da04aa4e998ace21b5cece135c7a794/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $addrMd3, "", "", "$newGuid.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), $addrMd3, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), $addrXlf3, "", "", "$newGuid.$newHash.de-de.xlf")

$wb.Save()
